$d = $word.ActiveDocument

$replacements = @(
    @("508×9=", "483×2="),
    @("202×9=", "965×4="),
    @("885×6=", "686×9="),
    @("630×4=", "196×8="),
    @("863×4=", "774×9="),
    @("832×8=", "896×2="),
    @("278×2=", "790×6="),
    @("480×9=", "483×2="),
    @("272×9=", "870×9="),
    @("338×9=", "898×6="),
    @("933×3=", "461×4="),
    @("742×2=", "447×2="),
    @("457×4=", "954×6="),
    @("834×2=", "102×2="),
    @("933×7=", "447×2="),
    @("845×9=", "985×9="),
    @("483×3=", "293×9="),
    @("437×4=", "657×5="),
    @("388×8=", "531×4="),
    @("423×9=", "774×2="),
    @("714×5=", "823×6="),
    @("848×4=", "479×9="),
    @("489×5=", "652×2="),
    @("102×3=", "136×8="),
    @("768×3=", "522×3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
